$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "AI"
$ws.Range("C2").Value = "test123"
$ws.Range("D2").Value = "23BCA001"
$ws.Range("E2").Value = "Anshika Bharti"
$ws.Range("F2").Value = "DEVICE_TEST"
$ws.Range("G2").Value = "127.0.0.1"
$ws.Range("H2").Value = "2025-07-03 12:58:10"
